# Update the "取得日時" (acquired datetime) column on the "ランサーズ" sheet.
# All data rows (2-18) currently hold the literal text "2025-09-09 01:15:32"
# in column A; they should be updated to "2025-09-09 01:44:37".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-09-09 01:15:32"
$newValue = "2025-09-09 01:44:37"

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = 18
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}

$wb.Save()
